$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5060
$ws.Range("J3").Value = 5380
$ws.Range("F4").Value = 1897
$ws.Range("G4").Value = 1470
$ws.Range("J4").Value = 1199
$ws.Range("J5").Value = 421
$ws.Range("J6").Value = 6700
$ws.Range("F7").Value = 24088
$ws.Range("G7").Value = 24694
$ws.Range("J7").Value = 18760

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 159
$ws.Range("J7").Value = 250

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 299
$ws.Range("J7").Value = 809

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 112
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 75
$ws.Range("J7").Value = 292

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 277
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 192
$ws.Range("J7").Value = 724

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 151
$ws.Range("J4").Value = 68
$ws.Range("G7").Value = 709
$ws.Range("J8").Value = 1197
$ws.Range("J11").Value = 288
$ws.Range("J14").Value = 88
$ws.Range("J15").Value = 201
$ws.Range("J19").Value = 536
$ws.Range("F20").Value = 645
$ws.Range("J24").Value = 57
$ws.Range("J25").Value = 93
$ws.Range("J29").Value = 1047
$ws.Range("J30").Value = 77
$ws.Range("J31").Value = 171
$ws.Range("J32").Value = 31
$ws.Range("J33").Value = 855
$ws.Range("J34").Value = 86
$ws.Range("J37").Value = 590
$ws.Range("J42").Value = 760
$ws.Range("J43").Value = 156
$ws.Range("J51").Value = 240
$ws.Range("J52").Value = 476
$ws.Range("J53").Value = 250
$ws.Range("J54").Value = 356
$ws.Range("J55").Value = 239
$ws.Range("J60").Value = 117
$ws.Range("J63").Value = 66
$ws.Range("J65").Value = 491
$ws.Range("J67").Value = 724
$ws.Range("J68").Value = 35
$ws.Range("J71").Value = 64
$ws.Range("J72").Value = 74
$ws.Range("J77").Value = 150
$ws.Range("J78").Value = 237
$ws.Range("J79").Value = 541
$ws.Range("J85").Value = 809
$ws.Range("J86").Value = 118
$ws.Range("J88").Value = 206
$ws.Range("J89").Value = 235
$ws.Range("J90").Value = 210
$ws.Range("J91").Value = 211
$ws.Range("J93").Value = 81
$ws.Range("J94").Value = 179
$ws.Range("J95").Value = 281
$ws.Range("J97").Value = 154
$ws.Range("J98").Value = 124
$ws.Range("J99").Value = 292
$ws.Range("F101").Value = 24088
$ws.Range("G101").Value = 24694
$ws.Range("J101").Value = 18760

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 179
$ws.Range("J4").Value = 16
$ws.Range("J7").Value = 590

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 213
$ws.Range("J7").Value = 855

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 144
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 174
$ws.Range("J7").Value = 491

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 69
$ws.Range("J6").Value = 166
$ws.Range("J7").Value = 356

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 312
$ws.Range("J3").Value = 363
$ws.Range("J4").Value = 59
$ws.Range("J6").Value = 272
$ws.Range("J7").Value = 1047

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 132
$ws.Range("J3").Value = 155
$ws.Range("J7").Value = 536

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("G4").Value = 40
$ws.Range("J4").Value = 21
$ws.Range("J6").Value = 174
$ws.Range("G7").Value = 709

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 151
$ws.Range("J6").Value = 390
$ws.Range("J7").Value = 760

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 64
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 211

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 192
$ws.Range("J4").Value = 34
$ws.Range("J7").Value = 541

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 133
$ws.Range("F4").Value = 49
$ws.Range("J4").Value = 39
$ws.Range("F7").Value = 645

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J3").Value = 29
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 111
$ws.Range("J3").Value = 140
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 476

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J3").Value = 18
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 63
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 210

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 66
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 41
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 156

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 336
$ws.Range("J3").Value = 360
$ws.Range("J6").Value = 395
$ws.Range("J7").Value = 1197

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 68
